$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the three "Expected Results" cells that referred to the old
#    hard-coded https://localhost:44302 URLs so they use the new relative,
#    auth-protected routes instead.
# ---------------------------------------------------------------------------
$ws.Cells.Item(4, 4).Value  = "User is redirected to /Quizzes/Index"
$ws.Cells.Item(5, 4).Value  = "User remains on /Account/Login"
$ws.Cells.Item(10, 4).Value = "User is redirected to /Account/Login"

# ---------------------------------------------------------------------------
# 2. Append four new scenario rows documenting that unauthenticated users get
#    redirected to the login page when navigating to the various Quizzes
#    routes. They are inserted directly below the existing last row (52) and
#    inherit that row's formatting (wrap text / column styles / row height).
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 4; $i++) {
    $ws.Rows.Item(53).Insert()
}

$ws.Rows.Item(52).Copy()
$ws.Range("A53:F56").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$newRows = @(
    @{ Row = 53; Num = 50; Scenario = "Unauthenticated user navigates to /Quizzes/Create" },
    @{ Row = 54; Num = 51; Scenario = "Unauthenticated user navigates to /Quizzes/Details?id=1" },
    @{ Row = 55; Num = 52; Scenario = "Unauthenticated user navigates to /Quizzes/Edit?id=1" },
    @{ Row = 56; Num = 53; Scenario = "Unauthenticated user navigates to /Quizzes/Delete?id=1" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Num
    $ws.Cells.Item($row, 2).Value = "Navigation"
    $ws.Cells.Item($row, 3).Value = $r.Scenario
    $ws.Cells.Item($row, 4).Value = "User is redirected to /Account/Login"
    $ws.Cells.Item($row, 5).Value = "Same as expected"
    $ws.Cells.Item($row, 6).Value = "Pass"
}

# ---------------------------------------------------------------------------
# 3. Extend the conditional formatting range that covered the last block of
#    rows (F49:F52) so it also covers the four newly added rows (F49:F56).
# ---------------------------------------------------------------------------
$lastBlock = $ws.Range("F49").FormatConditions
for ($i = 1; $i -le $lastBlock.Count(); $i++) {
    $lastBlock.Item($i).ModifyAppliesToRange($ws.Range("F49:F56"))
}

# ---------------------------------------------------------------------------
# 4. Update the view so the window shows the newly added rows, matching the
#    author's final cursor position/scroll state.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("E61").Select()
$excel.ActiveWindow.ScrollRow = 43

$wb.Save()
